# Insert a new data row at row 536 (pushes existing rows 536:618 down to 537:619)
# and populate it with a new Choclo price-report record, matching the rest of
# the "Vega Modelo de Temuco" sheet's constant columns (A,B,C,E,F,G,H,I,R).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(536).Insert()

$ws.Cells.Item(536, 1).Value = 10
$ws.Cells.Item(536, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(536, 3).Value = "La Araucanía"
$ws.Cells.Item(536, 4).Value = 45034
$ws.Cells.Item(536, 5).Value = 9
$ws.Cells.Item(536, 6).Value = 100112024
$ws.Cells.Item(536, 7).Value = "Choclo"
$ws.Cells.Item(536, 8).Value = "Dulce o Americano"
$ws.Cells.Item(536, 9).Value = "Primera"
$ws.Cells.Item(536, 10).Value = 5500
$ws.Cells.Item(536, 11).Value = 250
$ws.Cells.Item(536, 12).Value = 260
$ws.Cells.Item(536, 13).Value = 255
$ws.Cells.Item(536, 14).Value = "$/unidad"
$ws.Cells.Item(536, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(536, 16).Value = 255
$ws.Cells.Item(536, 17).Value = 1
$ws.Cells.Item(536, 18).Value = "Hortaliza"
